$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from E1 onto the new header cells F1:K1
$ws.Range("E1").Copy()
$ws.Range("F1:K1").PasteSpecial(-4122)

# Headers for new columns F1:K1
$ws.Range("F1").Value = "frequency"
$ws.Range("G1").Value = "frequency_occurrence"
$ws.Range("H1").Value = "frequency_occurrence_probab"
$ws.Range("I1").Value = "max_probab"
$ws.Range("J1").Value = "max_probab_percentage"
$ws.Range("K1").Value = "recommended_level"

# Data rows 2-38: F=frequency, G=frequency_occurrence, H=frequency_occurrence_probab,
# I=max_probab, J=max_probab_percentage (kept as text), K=recommended_level (moved from old F)
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = '{"L3":5,"L2":3}'
$ws.Range("H2").Value = '{"L3":0.625,"L2":0.375}'
$ws.Range("I2").Value = 0.625
$ws.Range("J2").Value = "'62.50"
$ws.Range("K2").Value = "L3"
$ws.Range("F3").Value = 16
$ws.Range("G3").Value = '{"L3":12,"L2":2,"L1":2}'
$ws.Range("H3").Value = '{"L3":0.75,"L2":0.125,"L1":0.125}'
$ws.Range("I3").Value = 0.75
$ws.Range("J3").Value = "'75.00"
$ws.Range("K3").Value = "L3"
$ws.Range("F4").Value = 8
$ws.Range("G4").Value = '{"L3":7,"L2":1}'
$ws.Range("H4").Value = '{"L3":0.875,"L2":0.125}'
$ws.Range("I4").Value = 0.875
$ws.Range("J4").Value = "'87.50"
$ws.Range("K4").Value = "L3"
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = '{"L2":5,"L1":2,"L3":1}'
$ws.Range("H5").Value = '{"L2":0.625,"L1":0.25,"L3":0.125}'
$ws.Range("I5").Value = 0.625
$ws.Range("J5").Value = "'62.50"
$ws.Range("K5").Value = "L3"
$ws.Range("F6").Value = 8
$ws.Range("G6").Value = '{"L3":8}'
$ws.Range("H6").Value = '{"L3":1.0}'
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = "'100.00"
$ws.Range("K6").Value = "L3"
$ws.Range("F7").Value = 8
$ws.Range("G7").Value = '{"L3":4,"L2":4}'
$ws.Range("H7").Value = '{"L3":0.5,"L2":0.5}'
$ws.Range("I7").Value = 0.5
$ws.Range("J7").Value = "'50.00"
$ws.Range("K7").Value = "L3"
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = '{"L3":5,"L2":3}'
$ws.Range("H8").Value = '{"L3":0.625,"L2":0.375}'
$ws.Range("I8").Value = 0.625
$ws.Range("J8").Value = "'62.50"
$ws.Range("K8").Value = "L3"
$ws.Range("F9").Value = 8
$ws.Range("G9").Value = '{"L3":5,"L2":3}'
$ws.Range("H9").Value = '{"L3":0.625,"L2":0.375}'
$ws.Range("I9").Value = 0.625
$ws.Range("J9").Value = "'62.50"
$ws.Range("K9").Value = "L3"
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = '{"L3":8}'
$ws.Range("H10").Value = '{"L3":1.0}'
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = "'100.00"
$ws.Range("K10").Value = "L3"
$ws.Range("F11").Value = 8
$ws.Range("G11").Value = '{"L3":7,"L2":1}'
$ws.Range("H11").Value = '{"L3":0.875,"L2":0.125}'
$ws.Range("I11").Value = 0.875
$ws.Range("J11").Value = "'87.50"
$ws.Range("K11").Value = "L3"
$ws.Range("F12").Value = 8
$ws.Range("G12").Value = '{"L2":5,"L3":3}'
$ws.Range("H12").Value = '{"L2":0.625,"L3":0.375}'
$ws.Range("I12").Value = 0.625
$ws.Range("J12").Value = "'62.50"
$ws.Range("K12").Value = "L3"
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = '{"L3":6,"L2":2}'
$ws.Range("H13").Value = '{"L3":0.75,"L2":0.25}'
$ws.Range("I13").Value = 0.75
$ws.Range("J13").Value = "'75.00"
$ws.Range("K13").Value = "L3"
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = '{"L2":5,"L3":3}'
$ws.Range("H14").Value = '{"L2":0.625,"L3":0.375}'
$ws.Range("I14").Value = 0.625
$ws.Range("J14").Value = "'62.50"
$ws.Range("K14").Value = "L3"
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = '{"L3":7,"L1":1}'
$ws.Range("H15").Value = '{"L3":0.875,"L1":0.125}'
$ws.Range("I15").Value = 0.875
$ws.Range("J15").Value = "'87.50"
$ws.Range("K15").Value = "L3"
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = '{"L3":6,"L2":2}'
$ws.Range("H16").Value = '{"L3":0.75,"L2":0.25}'
$ws.Range("I16").Value = 0.75
$ws.Range("J16").Value = "'75.00"
$ws.Range("K16").Value = "L3"
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = '{"L2":4,"L3":4}'
$ws.Range("H17").Value = '{"L2":0.5,"L3":0.5}'
$ws.Range("I17").Value = 0.5
$ws.Range("J17").Value = "'50.00"
$ws.Range("K17").Value = "L3"
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = '{"L3":8}'
$ws.Range("H18").Value = '{"L3":1.0}'
$ws.Range("I18").Value = 1
$ws.Range("J18").Value = "'100.00"
$ws.Range("K18").Value = "L3"
$ws.Range("F19").Value = 8
$ws.Range("G19").Value = '{"L3":8}'
$ws.Range("H19").Value = '{"L3":1.0}'
$ws.Range("I19").Value = 1
$ws.Range("J19").Value = "'100.00"
$ws.Range("K19").Value = "L3"
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = '{"L3":7,"L2":1}'
$ws.Range("H20").Value = '{"L3":0.875,"L2":0.125}'
$ws.Range("I20").Value = 0.875
$ws.Range("J20").Value = "'87.50"
$ws.Range("K20").Value = "L3"
$ws.Range("F21").Value = 8
$ws.Range("G21").Value = '{"L3":8}'
$ws.Range("H21").Value = '{"L3":1.0}'
$ws.Range("I21").Value = 1
$ws.Range("J21").Value = "'100.00"
$ws.Range("K21").Value = "L3"
$ws.Range("F22").Value = 8
$ws.Range("G22").Value = '{"L3":7,"L2":1}'
$ws.Range("H22").Value = '{"L3":0.875,"L2":0.125}'
$ws.Range("I22").Value = 0.875
$ws.Range("J22").Value = "'87.50"
$ws.Range("K22").Value = "L3"
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = '{"L3":5,"L2":3}'
$ws.Range("H23").Value = '{"L3":0.625,"L2":0.375}'
$ws.Range("I23").Value = 0.625
$ws.Range("J23").Value = "'62.50"
$ws.Range("K23").Value = "L3"
$ws.Range("F24").Value = 8
$ws.Range("G24").Value = '{"L3":8}'
$ws.Range("H24").Value = '{"L3":1.0}'
$ws.Range("I24").Value = 1
$ws.Range("J24").Value = "'100.00"
$ws.Range("K24").Value = "L3"
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = '{"L3":5,"L2":3}'
$ws.Range("H25").Value = '{"L3":0.625,"L2":0.375}'
$ws.Range("I25").Value = 0.625
$ws.Range("J25").Value = "'62.50"
$ws.Range("K25").Value = "L3"
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = '{"L3":8}'
$ws.Range("H26").Value = '{"L3":1.0}'
$ws.Range("I26").Value = 1
$ws.Range("J26").Value = "'100.00"
$ws.Range("K26").Value = "L3"
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = '{"L2":5,"L3":3}'
$ws.Range("H27").Value = '{"L2":0.625,"L3":0.375}'
$ws.Range("I27").Value = 0.625
$ws.Range("J27").Value = "'62.50"
$ws.Range("K27").Value = "L3"
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = '{"L3":4,"L2":3,"L1":1}'
$ws.Range("H28").Value = '{"L3":0.5,"L2":0.375,"L1":0.125}'
$ws.Range("I28").Value = 0.5
$ws.Range("J28").Value = "'50.00"
$ws.Range("K28").Value = "L3"
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = '{"L2":8}'
$ws.Range("H29").Value = '{"L2":1.0}'
$ws.Range("I29").Value = 1
$ws.Range("J29").Value = "'100.00"
$ws.Range("K29").Value = "L2"
$ws.Range("F30").Value = 8
$ws.Range("G30").Value = '{"L3":8}'
$ws.Range("H30").Value = '{"L3":1.0}'
$ws.Range("I30").Value = 1
$ws.Range("J30").Value = "'100.00"
$ws.Range("K30").Value = "L3"
$ws.Range("F31").Value = 8
$ws.Range("G31").Value = '{"L3":8}'
$ws.Range("H31").Value = '{"L3":1.0}'
$ws.Range("I31").Value = 1
$ws.Range("J31").Value = "'100.00"
$ws.Range("K31").Value = "L3"
$ws.Range("F32").Value = 8
$ws.Range("G32").Value = '{"L3":8}'
$ws.Range("H32").Value = '{"L3":1.0}'
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = "'100.00"
$ws.Range("K32").Value = "L3"
$ws.Range("F33").Value = 8
$ws.Range("G33").Value = '{"L3":8}'
$ws.Range("H33").Value = '{"L3":1.0}'
$ws.Range("I33").Value = 1
$ws.Range("J33").Value = "'100.00"
$ws.Range("K33").Value = "L3"
$ws.Range("F34").Value = 8
$ws.Range("G34").Value = '{"L2":5,"L3":3}'
$ws.Range("H34").Value = '{"L2":0.625,"L3":0.375}'
$ws.Range("I34").Value = 0.625
$ws.Range("J34").Value = "'62.50"
$ws.Range("K34").Value = "L3"
$ws.Range("F35").Value = 8
$ws.Range("G35").Value = '{"L3":8}'
$ws.Range("H35").Value = '{"L3":1.0}'
$ws.Range("I35").Value = 1
$ws.Range("J35").Value = "'100.00"
$ws.Range("K35").Value = "L3"
$ws.Range("F36").Value = 8
$ws.Range("G36").Value = '{"L3":8}'
$ws.Range("H36").Value = '{"L3":1.0}'
$ws.Range("I36").Value = 1
$ws.Range("J36").Value = "'100.00"
$ws.Range("K36").Value = "L3"
$ws.Range("F37").Value = 8
$ws.Range("G37").Value = '{"L3":8}'
$ws.Range("H37").Value = '{"L3":1.0}'
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = "'100.00"
$ws.Range("K37").Value = "L3"
$ws.Range("F38").Value = 8
$ws.Range("G38").Value = '{"L3":8}'
$ws.Range("H38").Value = '{"L3":1.0}'
$ws.Range("I38").Value = 1
$ws.Range("J38").Value = "'100.00"
$ws.Range("K38").Value = "L3"
